$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows (product_id, rate, scope) replacing/extending the existing
# price-list starting at row 2.
$rows = @(
    @("Blood",      112, "All"),
    @("Clementine",  113, "All"),
    @("Grapefruit",  88,  "All"),
    @("Mandarin",    102, "43"),
    @("Mandarin",    120, "45"),
    @("Mandarin",    104, "All"),
    @("Navel",       93,  "All"),
    @("Shamuti",     84,  "All"),
    @("Tangerine",   85,  "12"),
    @("Tangerine",   92,  "All"),
    @("Valencia",    90,  "45"),
    @("Valencia",    87,  "All")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    $scopeCell = $ws.Cells.Item($r, 3)
    $scope = $row[2]
    # Scope values that look numeric ("43", "45", "12", ...) must stay text,
    # matching the source data ("All" is already text, so format as text
    # before assignment to stop Excel's auto-number inference).
    if ($scope -match '^[0-9]+$') {
        $scopeCell.NumberFormat = "@"
    }
    $scopeCell.Value = $scope

    $r++
}
